$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 5+ down by one so a blank row appears at row 5 and the
# jx:if command row (with its highlight style) lands on row 6 -------------
$ws.Rows.Item(5).Insert()

# --- Text content ----------------------------------------------------------
# Row 1: title (unchanged)
$ws.Range("A1").Value = "Person Report"

# Row 3: Name label/value split across two columns
$ws.Range("A3").Value = "Name:"
$ws.Range("B3").Value = '${person.name}'

# Row 4: Age label/value split across two columns
$ws.Range("A4").Value = "Age:"
$ws.Range("B4").Value = '${person.age}'

# Row 6 (was row 5): jx:if command, new argument syntax
$ws.Range("A6").Value = 'jx:if(condition="person.age < 18", lastCell="B6")'

# Row 7 (was row 6): Parent label/value split across two columns
$ws.Range("A7").Value = "Parent:"
$ws.Range("B7").Value = '${person.parentName}'

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 22.666666666666668

# --- Highlight fill colour on the jx:if row (indexed 22 -> 43) --------
$ws.Range("A6").Interior.ColorIndex = 43

# --- Explanatory cell comment on the jx:if row --------------------------
$ws.Range("A6").AddComment("JXLS Command:`nThis row contains the jx:if condition.`nIf condition is true, the next row(s) will be included.`nIf false, they will be removed from output.")
